$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ajoute des variables de session (responsables) dans la page de login, admin, vet et rec
$ws.Range("B21").Value = "Fosso"
$ws.Range("B25").Value = "Alejandro"
$ws.Range("B26").Value = "Alejandro"

# Met a jour la selection active
$ws.Range("C25").Select()
